# Update cryptos list data (Price column D, Volume(1h) column E) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.287.81"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.932.99"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.68"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.01"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("E8").Value = "  -0.36%  "

$ws.Range("E9").Value = "  +3.13%  "

$ws.Range("E10").Value = "  -0.55%  "

$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("E12").Value = "  -0.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.64"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.04%  "

$ws.Range("E14").Value = "  +0.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.418.03"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.201.81"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.74"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.933.59"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "433.18"
$ws.Range("D19").ClearFormats()

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.51"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.681"
$ws.Range("D21").ClearFormats()

$ws.Range("E22").Value = "  +0.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "82.01"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.01"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.49%  "

$ws.Range("E25").Value = "  -0.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.88"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.02%  "

$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("E28").Value = "  -4.18%  "

$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.01"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.19%  "

$ws.Range("E31").Value = "  +3.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.80"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.98%  "

$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0890"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.49%  "

$ws.Range("E35").Value = "  -0.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.68"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.01"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.01"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.34%  "

$ws.Range("E39").Value = "  +0.06%  "

$ws.Range("E40").Value = "  +0.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.44"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +7.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.284"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0349"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "371.49"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.702.66"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.60"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.28%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.94"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.96%  "

$ws.Range("E49").Value = "  -0.99%  "

$ws.Range("E50").Value = "  -0.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.126"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.26%  "
